$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted ahead of the existing row 513
# (Macroferia Regional de Talca - Zanahoria), pushing rows 513:530 down to
# 514:531. Insert a physical row at 513 so everything below shifts down,
# then populate the new row with the same fixed descriptive data as the
# rest of the block and the new record's own date / price figures.
$ws.Rows.Item(513).Insert()

$ws.Range("A513").Value = 5
$ws.Range("B513").Value = "Macroferia Regional de Talca"
$ws.Range("C513").Value = "Maule"
$ws.Range("D513").Value2 = 45075
$ws.Range("E513").Value = 7
$ws.Range("F513").Value = 100114013
$ws.Range("G513").Value = "Zanahoria"
$ws.Range("H513").Value = "Sin especificar"
$ws.Range("I513").Value = "Primera"
$ws.Range("J513").Value = 500
$ws.Range("K513").Value = 5000
$ws.Range("L513").Value = 5000
$ws.Range("M513").Value = 5000
$ws.Range("N513").Value = "$/saco 20 kilos"
$ws.Range("O513").Value = "Región de Ñuble"
$ws.Range("P513").Value = 250
$ws.Range("Q513").Value = 20
$ws.Range("R513").Value = "Hortaliza"
